$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1578, 1578, 1578, 0.01048845450083415),
    @(1, 911, 911, 911, 0.01050256093343099),
    @(2, 1233, 1233, 1233, 0.01222229798634847),
    @(3, 1162, 1162, 1162, 0.0119592030843099),
    @(4, 1356, 1356, 1356, 0.01250503857930501),
    @(5, 1484, 1484, 1484, 0.01210121313730876),
    @(6, 1921, 1921, 1921, 0.01204830805460612),
    @(7, 1748, 1748, 1748, 0.010457444190979),
    @(8, 1645, 1645, 1645, 0.01208465894063314),
    @(9, 1741, 1741, 1741, 0.01212320327758789)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
